# Fix missed remaining label
#
# Adds the "remaining" textbox that was missing from the activity diagram
# (mirrors the sibling "commandWord" textbox already on the slide).
#
# EMU targets (from the authoritative OOXML):
#   off  x = 6102496   y = -20104572
#   ext cx = 1972764  cy = 369460
# PowerPoint's Shape.Left/Top/Width/Height are points, Single-precision;
# feeding them straight into AddTextbox (rather than re-assigning the
# properties afterwards) keeps the EMU round-trip exact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$left   = 6102496 / 12700
$top    = -20104572 / 12700
$width  = 1972764 / 12700
$height = 369460 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 125"

$tr = $tb.TextFrame.TextRange
$tr.Text = "remaining"
$tr.Font.Size = 18.01
$tr.LanguageID = "en-SG"

$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

$tb.Fill.Visible = $false

# AutoSize re-derives Width/Height from the text metrics; pin them back to
# the exact authored extent (Left/Top are already exact from AddTextbox,
# so they are intentionally left untouched).
$tb.Width = $width
$tb.Height = $height
